$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Backlog")
$ws.Activate() | Out-Null

# Row 11: "Tax Accountants on SATURDAY 9/8" (IN PROGRESS) -> "Write check for $3600 to State of IL" (TODO)
$ws.Range("A11").Value = "Write check for `$3600 to State of IL"
$ws.Range("B11").Value = "TODO"

# Row 14's status becomes "TUESDAY 9/11" (new shared string, written here so it lands right
# after the "Write check..." string in the rebuilt shared-string table).
$ws.Range("B14").Value = "TUESDAY 9/11"

# Row 12: "Bank of America Maintenance Fees" -> new row "Download and start using YNAB" (TODO)
$ws.Range("A12").Value = "Download and start using YNAB"
$ws.Range("B12").Value = "TODO"

# Row 13: "Goto the Doctor and get checked up - Call Kachar" -> "Bank of America Maintenance Fees" (TODO)
$ws.Range("A13").Value = "Bank of America Maintenance Fees"
$ws.Range("B13").Value = "TODO"

# Row 14: "House Hunting Backlog" (special highlighted row) -> "Goto the Doctor and get checked up - Call Kachar",
# restyled to match the normal yellow-highlight rows (same look as rows 11-13/18).
$ws.Range("A14").Value = "Goto the Doctor and get checked up - Call Kachar"

$ws.Range("A14").Interior.Color = 65535
$ws.Range("A14").Font.Bold = $true
$ws.Range("B14").Interior.Color = 65535
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").HorizontalAlignment = -4131

# Leave the selection where the user last clicked while editing.
$ws.Range("B12").Select() | Out-Null
